# Updates 'want-to-go' counts (F) and minimum-price values (G)
# across the '展览', '演出' and '全部类型' sheets,
# matching the refreshed data snapshot output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 824
$ws.Range("G2").Value = 58
$ws.Range("F5").Value = 1162
$ws.Range("F6").Value = 3386
$ws.Range("F7").Value = 2491
$ws.Range("F8").Value = 63
$ws.Range("F9").Value = 2320
$ws.Range("F10").Value = 244
$ws.Range("F13").Value = 1617
$ws.Range("F14").Value = 625
$ws.Range("F15").Value = 86
$ws.Range("F16").Value = 290
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 42
$ws.Range("F20").Value = 61
$ws.Range("F21").Value = 418
$ws.Range("F22").Value = 20
$ws.Range("F24").Value = 457
$ws.Range("F25").Value = 650
$ws.Range("F26").Value = 71
$ws.Range("F27").Value = 66
$ws.Range("F28").Value = 345
$ws.Range("F29").Value = 29
$ws.Range("F30").Value = 1601
$ws.Range("F31").Value = 764
$ws.Range("F32").Value = 774
$ws.Range("F33").Value = 1890
$ws.Range("F34").Value = 207
$ws.Range("F35").Value = 491
$ws.Range("F36").Value = 71
$ws.Range("F37").Value = 538
$ws.Range("F38").Value = 1174

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 56
$ws.Range("F7").Value = 9

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 824
$ws.Range("G2").Value = 58
$ws.Range("F5").Value = 1162
$ws.Range("F6").Value = 3386
$ws.Range("F7").Value = 2491
$ws.Range("F8").Value = 63
$ws.Range("F9").Value = 2320
$ws.Range("F10").Value = 244
$ws.Range("F13").Value = 1617
$ws.Range("F14").Value = 625
$ws.Range("F15").Value = 86
$ws.Range("F16").Value = 290
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 42
$ws.Range("F20").Value = 61
$ws.Range("F21").Value = 418
$ws.Range("F22").Value = 20
$ws.Range("F24").Value = 457
$ws.Range("F25").Value = 650
$ws.Range("F26").Value = 71
$ws.Range("F27").Value = 56
$ws.Range("F30").Value = 66
$ws.Range("F31").Value = 345
$ws.Range("F32").Value = 29
$ws.Range("F33").Value = 1601
$ws.Range("F34").Value = 764
$ws.Range("F36").Value = 774
$ws.Range("F37").Value = 1890
$ws.Range("F38").Value = 207
$ws.Range("F40").Value = 9
$ws.Range("F42").Value = 491
$ws.Range("F43").Value = 71
$ws.Range("F44").Value = 538
$ws.Range("F45").Value = 1174

